# News From The Forest Links: add the April 2019 newsletter row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aprilUrl = "https://myemail.constantcontact.com/News-From-The-Forest--April-2019.html?soid=1102494320279&aid=oK_3DP6m1cU"

# New row 9: month label in A, newsletter link in B (mirrors existing rows).
$ws.Range("A9").Value = "April 2019"
$ws.Range("B9").Value = $aprilUrl

# Turn B9 into a real hyperlink, then restyle it to match the other link
# cells (the Add() call stamps its own style variant onto the cell first).
$ws.Hyperlinks.Add($ws.Range("B9"), $aprilUrl)
$ws.Range("B9").Style = $ws.Range("B7").Style

# Matches the author's final selection position in the saved file.
$ws.Range("B19").Select()
